# Adapt tests to control version
# The user adds a "version" column to the settings sheet (header "version"
# in C1, value "e" in C2) and leaves the "settings" sheet as the active
# sheet/selection sitting on C3 (just past the newly typed data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

$ws.Activate()

$ws.Range("C1").Value = "version"
$ws.Range("C2").Value = "e"

$ws.Range("C3").Select()
